# Generate Report for Handback
#
# The handback step completed: the translated files (xlf) came back "in
# sync with en-US", so:
#   1. Update the Status text from "Ready for handoff" to
#      "Handed back: in sync with en-US" everywhere it is shown
#      (Overview sheet + per-language detail sheets).
#   2. Record the newly produced "Latest Target File" / "Latest Handback
#      File" for each row on the per-language detail sheets, with the
#      same hyperlinked-file look the existing File columns use.
#   3. Stamp the real "Latest Handback DateTime" for each row (previously
#      the zero/placeholder date), different per language since they
#      completed the handback pipeline at different times.

$wb = $excel.ActiveWorkbook
$ov = $wb.Worksheets.Item(1)   # Overview
$zh = $wb.Worksheets.Item(2)   # zh-cn
$de = $wb.Worksheets.Item(3)   # de-de

$newStatus = "Handed back: in sync with en-US"

# --- 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# --- 2. zh-cn: populate "Latest Target File" (F) / "Latest Handback File" (G) ---
$zhMdName  = "9c8d3236-453f-43c2-b43e-551b6fffcf33.md"
$zhMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/f74873dc56204db750e19f38a9233181ce15d4d9/e2e/9c8d3236-453f-43c2-b43e-551b6fffcf33.md"
$zhXlfName = "9c8d3236-453f-43c2-b43e-551b6fffcf33.efea525a0597b57715c41c3d624c6815f4451374.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78732aa4cb00d7c9ae2e8034c82a35d53ca702ee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9c8d3236-453f-43c2-b43e-551b6fffcf33.efea525a0597b57715c41c3d624c6815f4451374.zh-cn.xlf"

foreach ($row in 2,3) {
    $fCell = $zh.Range("F$row")
    $fCell.Value = $zhMdName
    $zh.Hyperlinks.Add($fCell, $zhMdUrl, [Type]::Missing, [Type]::Missing, $zhMdName)
    $fCell.Style = "HyperLink"

    $gCell = $zh.Range("G$row")
    $gCell.Value = $zhXlfName
    $zh.Hyperlinks.Add($gCell, $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)
    $gCell.Style = "HyperLink"
}

# zh-cn: stamp the real "Latest Handback DateTime" (was the 0001-01-01 placeholder)
$zh.Range("H2").Value = "2016-03-22 11:07:43"
$zh.Range("H3").Value = "2016-03-22 11:07:43"

# --- 3. de-de: populate "Latest Target File" (F) / "Latest Handback File" (G) ---
$deMdName  = "9c8d3236-453f-43c2-b43e-551b6fffcf33.md"
$deMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/f74873dc56204db750e19f38a9233181ce15d4d9/e2e/9c8d3236-453f-43c2-b43e-551b6fffcf33.md"
$deXlfName = "9c8d3236-453f-43c2-b43e-551b6fffcf33.efea525a0597b57715c41c3d624c6815f4451374.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/376082c6c69691d613a2c21445e263b3a9733d60/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9c8d3236-453f-43c2-b43e-551b6fffcf33.efea525a0597b57715c41c3d624c6815f4451374.de-de.xlf"

foreach ($row in 2,3) {
    $fCell = $de.Range("F$row")
    $fCell.Value = $deMdName
    $de.Hyperlinks.Add($fCell, $deMdUrl, [Type]::Missing, [Type]::Missing, $deMdName)
    $fCell.Style = "HyperLink"

    $gCell = $de.Range("G$row")
    $gCell.Value = $deXlfName
    $de.Hyperlinks.Add($gCell, $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfName)
    $gCell.Style = "HyperLink"
}

# de-de: stamp the real "Latest Handback DateTime" (later than zh-cn's)
$de.Range("H2").Value = "2016-03-22 11:07:50"
$de.Range("H3").Value = "2016-03-22 11:07:50"
